$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E4:S4").Value = "-"
$ws.Range("E4:S4").HorizontalAlignment = -4152

$ws.Range("E53:S53").Value = "-"
$ws.Range("E53:S53").HorizontalAlignment = -4152

$ws.Range("E93:S93").Value = "-"
$ws.Range("E93:S93").HorizontalAlignment = -4152

$ws.Range("J69").Value = "-"
$ws.Range("J69").HorizontalAlignment = -4152

$ws.Range("P69").Value = "-"
$ws.Range("P69").HorizontalAlignment = -4152

$ws.Range("S69").Value = "-"
$ws.Range("S69").HorizontalAlignment = -4152
